$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1089.8572
$ws.Range("I101").Value = 395
$ws.Range("J101").Value = 2016.3334
$ws.Range("K101").Value = 1185
$ws.Range("L101").Value = 6049.0002
$ws.Range("M101").Value = 437
$ws.Range("N101").Value = -9293.0002

$ws.Range("H107").Value = 869.43243
$ws.Range("I107").Value = 865.129
$ws.Range("J107").Value = 891.6667
$ws.Range("K107").Value = 865.129
$ws.Range("L107").Value = 891.6667
$ws.Range("M107").Value = 1054.871
$ws.Range("N107").Value = -4731.6667

$ws.Range("H113").Value = 1927.2727
$ws.Range("I113").Value = 2260
$ws.Range("J113").Value = 1650
$ws.Range("K113").Value = 2260
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 994
$ws.Range("N113").Value = -8158

$ws.Range("H120").Value = 46630
$ws.Range("J120").Value = 46630
$ws.Range("L120").Value = 46630
$ws.Range("N120").Value = -56306

$ws.Range("H127").Value = 1230
$ws.Range("J127").Value = 1352.9412
$ws.Range("L127").Value = 4058.8236
$ws.Range("N127").Value = -13978.8236

$ws.Range("H129").Value = 765.5
$ws.Range("J129").Value = 906
$ws.Range("L129").Value = 2718
$ws.Range("N129").Value = -12718

$ws.Range("H135").Value = 15146764
$ws.Range("I135").Value = 3807.5386
$ws.Range("J135").Value = 58893084
$ws.Range("K135").Value = 34267.8474
$ws.Range("L135").Value = 530037756
$ws.Range("M135").Value = -31732.8474
$ws.Range("N135").Value = -530042826

$ws.Range("H138").Value = 3231
$ws.Range("I138").Value = 2658.3333
$ws.Range("J138").Value = 3803.6667
$ws.Range("K138").Value = 7974.999899999999
$ws.Range("L138").Value = 11411.0001
$ws.Range("M138").Value = -2834.999899999999
$ws.Range("N138").Value = -21691.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 106497
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H37").Value = 9310.9
$ws.Range("J37").Value = 10119.444
$ws.Range("L37").Value = 10119.444
$ws.Range("N37").Value = -10665.444

$ws.Range("H45").Value = 1115.12
$ws.Range("I45").Value = 929.9091
$ws.Range("J45").Value = 1260.6428
$ws.Range("K45").Value = 929.9091
$ws.Range("L45").Value = 1260.6428
$ws.Range("M45").Value = -552.9091
$ws.Range("N45").Value = -2014.6428

$ws.Range("H68").Value = 62025.332
$ws.Range("J68").Value = 62025.332
$ws.Range("L68").Value = 62025.332
$ws.Range("N68").Value = -63647.332

$ws.Range("H71").Value = 62025.332
$ws.Range("J71").Value = 62025.332
$ws.Range("L71").Value = 186075.996
$ws.Range("N71").Value = -194187.996

$ws.Range("H74").Value = 15630598
$ws.Range("I74").Value = 22728014
$ws.Range("J74").Value = 16281.4
$ws.Range("K74").Value = 22728014
$ws.Range("L74").Value = 16281.4
$ws.Range("M74").Value = -22727140
$ws.Range("N74").Value = -18029.4

$ws.Range("H77").Value = 15630598
$ws.Range("I77").Value = 22728014
$ws.Range("J77").Value = 16281.4
$ws.Range("K77").Value = 113640070
$ws.Range("L77").Value = 81407
$ws.Range("M77").Value = -113635702
$ws.Range("N77").Value = -90143

$ws.Range("H139").Value = 70747
$ws.Range("J139").Value = 70747
$ws.Range("L139").Value = 70747
$ws.Range("N139").Value = -81027

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1798.7693
$ws.Range("I86").Value = 1668.1052
$ws.Range("J86").Value = 2153.4285
$ws.Range("K86").Value = 1668.1052
$ws.Range("L86").Value = 2153.4285
$ws.Range("M86").Value = -545.1052
$ws.Range("N86").Value = -4399.4285

$ws.Range("H89").Value = 1798.7693
$ws.Range("I89").Value = 1668.1052
$ws.Range("J89").Value = 2153.4285
$ws.Range("K89").Value = 8340.526
$ws.Range("L89").Value = 10767.1425
$ws.Range("M89").Value = -2724.526
$ws.Range("N89").Value = -21999.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9779.1
$ws.Range("J51").Value = 10023.875
$ws.Range("L51").Value = 10023.875
$ws.Range("N51").Value = -11495.875

$ws.Range("H58").Value = 2439.472
$ws.Range("I58").Value = 877.9143
$ws.Range("J58").Value = 3451.5925
$ws.Range("K58").Value = 877.9143
$ws.Range("L58").Value = 3451.5925
$ws.Range("M58").Value = -674.9143
$ws.Range("N58").Value = -3857.5925

$ws.Range("H61").Value = 9779.1
$ws.Range("J61").Value = 10023.875
$ws.Range("L61").Value = 10023.875
$ws.Range("N61").Value = -10719.875

$ws.Range("H68").Value = 18032.3
$ws.Range("J68").Value = 18450.555
$ws.Range("L68").Value = 18450.555
$ws.Range("N68").Value = -19948.555

$ws.Range("H71").Value = 18032.3
$ws.Range("J71").Value = 18450.555
$ws.Range("L71").Value = 55351.665
$ws.Range("N71").Value = -62839.665

$ws.Range("H74").Value = 14418
$ws.Range("J74").Value = 16469.777
$ws.Range("L74").Value = 16469.777
$ws.Range("N74").Value = -18217.777

$ws.Range("H77").Value = 14418
$ws.Range("J77").Value = 16469.777
$ws.Range("L77").Value = 49409.33099999999
$ws.Range("N77").Value = -58145.33099999999

$ws.Range("H94").Value = 2038.25
$ws.Range("I94").Value = 1112
$ws.Range("K94").Value = 1112
$ws.Range("M94").Value = -661

$ws.Range("H136").Value = 2439.472
$ws.Range("I136").Value = 877.9143
$ws.Range("J136").Value = 3451.5925
$ws.Range("K136").Value = 2633.7429
$ws.Range("L136").Value = 10354.7775
$ws.Range("M136").Value = -83.74290000000019
$ws.Range("N136").Value = -15454.7775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 861.41174
$ws.Range("I68").Value = 824.5714
$ws.Range("J68").Value = 887.2
$ws.Range("K68").Value = 2473.7142
$ws.Range("L68").Value = 2661.6
$ws.Range("M68").Value = -1662.7142
$ws.Range("N68").Value = -4283.6

$ws.Range("H71").Value = 861.41174
$ws.Range("I71").Value = 824.5714
$ws.Range("J71").Value = 887.2
$ws.Range("K71").Value = 7421.1426
$ws.Range("L71").Value = 7984.8
$ws.Range("M71").Value = -3365.1426
$ws.Range("N71").Value = -16096.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125222.89
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 222801.2
$ws.Range("K80").Value = 3250
$ws.Range("L80").Value = 222801.2
$ws.Range("M80").Value = -2252
$ws.Range("N80").Value = -224797.2

$ws.Range("H83").Value = 125222.89
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 222801.2
$ws.Range("K83").Value = 16250
$ws.Range("L83").Value = 1114006
$ws.Range("M83").Value = -11258
$ws.Range("N83").Value = -1123990

$ws.Range("H126").Value = 2000.1538
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 2250.25
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 6750.75
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -11690.75

$ws.Range("H132").Value = 2210.889
$ws.Range("I132").Value = 1316.6666
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 3949.9998
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -1419.9998
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2090.3333
$ws.Range("I7").Value = 2115.4285
$ws.Range("J7").Value = 2002.5
$ws.Range("K7").Value = 2115.4285
$ws.Range("L7").Value = 2002.5
$ws.Range("M7").Value = -2003.4285
$ws.Range("N7").Value = -2226.5

$ws.Range("H40").Value = 2625.35
$ws.Range("I40").Value = 2431.2
$ws.Range("J40").Value = 2819.5
$ws.Range("K40").Value = 2431.2
$ws.Range("L40").Value = 2819.5
$ws.Range("M40").Value = -2295.2
$ws.Range("N40").Value = -3091.5

$ws.Range("H100").Value = 2220
$ws.Range("I100").Value = 2120
$ws.Range("J100").Value = 2345
$ws.Range("K100").Value = 2120
$ws.Range("L100").Value = 2345
$ws.Range("M100").Value = -1579
$ws.Range("N100").Value = -3427

$ws.Range("H122").Value = 1972.875
$ws.Range("I122").Value = 1783.1072
$ws.Range("K122").Value = 5349.321599999999
$ws.Range("M122").Value = -2899.321599999999

$ws.Range("H126").Value = 2090.3333
$ws.Range("I126").Value = 2115.4285
$ws.Range("J126").Value = 2002.5
$ws.Range("K126").Value = 6346.2855
$ws.Range("L126").Value = 6007.5
$ws.Range("M126").Value = -3876.2855
$ws.Range("N126").Value = -10947.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 615.0909
$ws.Range("I107").Value = 435.1111
$ws.Range("K107").Value = 1305.3333
$ws.Range("M107").Value = 614.6667
